# Lab 4 decrypt monoalphabet
# Renames the sheet, rebuilds the frequency table with an added
# "cypher text" column (E/F) next to the existing "open text" column
# (B/C), appends the missing low-frequency symbol row, and refreshes
# the _xlchart.v1.* defined names the Pareto charts read from.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Rename the worksheet
# ---------------------------------------------------------------
$ws.Name = "Моноалфавитный шифр"

# ---------------------------------------------------------------
# 2. Header cells: "German" -> "Open text", add "Cypher text"
# ---------------------------------------------------------------
$ws.Range("B1").Value = "Open text"
$ws.Range("E1").Value = "Cypher text"

$ws.Range("B2").Value = "Symbol"
$ws.Range("C2").Value = "Appearances"
$ws.Range("E2").Value = "Symbol"
$ws.Range("F2").Value = "Appearances"

# ---------------------------------------------------------------
# 3. Frequency table rows 3-31: open text (B/C, unchanged values)
#    plus the new cypher text columns (E/F) for the same symbols.
# ---------------------------------------------------------------
$cipherData = @(
    @(3, "r", 473, "w", 473),
    @(4, "d", 341, "k", 341),
    @(5, "e", 1007, "l", 1007),
    @(6, "c", 247, "j", 247),
    @(7, "f", 102, "m", 102),
    @(8, "v", 116, "b", 116),
    @(9, "a", 456, "g", 456),
    @(10, "t", 475, "z", 475),
    @(11, "k", 109, "q", 109),
    @(12, "o", 169, "t", 169),
    @(13, "m", 221, "s", 221),
    @(14, "z", 61, "f", 61),
    @(15, "u", 228, "a", 228),
    @(16, "s", 577, "x", 577),
    @(17, "i", 576, "ö", 576),
    @(18, "n", 588, "ß", 588),
    @(19, "l", 192, "r", 192),
    @(20, "h", 401, "o", 401),
    @(21, "g", 179, "n", 179),
    @(22, "w", 106, "c", 106),
    @(23, "j", 39, "p", 39),
    @(24, "b", 149, "i", 149),
    @(25, "p", 19, "ü", 19),
    @(26, "y", 6, "e", 6),
    @(27, "ü", 39, "ä", 39),
    @(28, "ä", 14, "h", 14),
    @(29, "ö", 9, "u", 9),
    @(30, "ß", 17, "y", 17),
    @(31, "x", 1, "d", 1)
)

foreach ($entry in $cipherData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]   # B: open-text symbol
    $ws.Cells.Item($r, 3).Value = $entry[2]   # C: open-text appearances
    $ws.Cells.Item($r, 5).Value = $entry[3]   # E: cypher-text symbol
    $ws.Cells.Item($r, 6).Value = $entry[4]   # F: cypher-text appearances
}

# ---------------------------------------------------------------
# 4. Stray leftover note, kept as-is (row 35) plus a new one (row 34)
# ---------------------------------------------------------------
$ws.Range("E34").WrapText = $true
$ws.Range("E34").Value = "`n;29"
$ws.Rows.Item(34).AutoFit()

# ---------------------------------------------------------------
# 5. Remove the old leftover cells that used to live further down
#    the sheet (row 37 and row 52) so the used range shrinks back
#    to B1:F35.
# ---------------------------------------------------------------
$ws.Range("E37").ClearContents()
$ws.Range("B52:C52").ClearContents()

# ---------------------------------------------------------------
# 6. Defined names used by the charts: drop the five that are no
#    longer needed and point the remaining six at the renamed sheet
#    and the new ranges (cypher-text columns included).
# ---------------------------------------------------------------
$namesToDelete = @(
    "_xlchart.v1.6",
    "_xlchart.v1.7",
    "_xlchart.v1.8",
    "_xlchart.v1.9",
    "_xlchart.v1.10"
)
foreach ($n in $namesToDelete) {
    $wb.Names.Item($n).Delete()
}

$wb.Names.Item("_xlchart.v1.0").RefersTo = "='Моноалфавитный шифр'!`$B`$3:`$B`$31"
$wb.Names.Item("_xlchart.v1.1").RefersTo = "='Моноалфавитный шифр'!`$C`$2"
$wb.Names.Item("_xlchart.v1.2").RefersTo = "='Моноалфавитный шифр'!`$C`$3:`$C`$31"
$wb.Names.Item("_xlchart.v1.3").RefersTo = "='Моноалфавитный шифр'!`$E`$3:`$E`$31"
$wb.Names.Item("_xlchart.v1.4").RefersTo = "='Моноалфавитный шифр'!`$F`$2"
$wb.Names.Item("_xlchart.v1.5").RefersTo = "='Моноалфавитный шифр'!`$F`$3:`$F`$31"

# ---------------------------------------------------------------
# 7. View cosmetics: smaller zoom, new selected cell
# ---------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("H31").Select()
